# Update profit files after running on 2025-11-07
# Append a new row (row 67) with the latest allocation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67

# Column A holds a plain text date label (matches the style of existing rows,
# which store dates as literal text rather than Excel date serials), so we
# force the cell to a text format before assigning the value and then drop
# the number-format override again so no extra style is left behind.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "11/07/2025"
$ws.Range("A" + $newRow).Style = "Normal"

$ws.Range("B" + $newRow).Value = 0.1940232686290251
$ws.Range("C" + $newRow).Value = 0.8059767313709749
